# The closing paragraph hard-codes a specific issue date ("a los (24) dias
# del mes de (enero) de 2024."). Turn the day, month and year into
# template placeholders so the certificate can be generated dynamically,
# the same way the rest of the letter already uses ${name}, ${post}, etc.
#
# Word COM Find.Execute signature:
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

$d = $word.ActiveDocument

# "(24)" -> "(${day})"
$d.Content.Find.Execute("(24)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(`${day})", 2)

# "(enero)" -> "(${month})"
$d.Content.Find.Execute("(enero)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(`${month})", 2)

# "de 2024." -> "de ${year}." (collapses the old literal year + trailing period)
$d.Content.Find.Execute("de 2024.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "de `${year}.", 2)
